$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-parsed as a number by Excel
# are first marked as Text (NumberFormat "@") so they are stored as strings,
# matching the workbook-wide convention of storing these figures as text.

$ws.Range("D2").Value = "27.402.90"
$ws.Range("E2").Value = "  -1.26%  "
$ws.Range("D3").Value = "1.711.12"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.42"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("E6").Value = "  -2.25%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2670"
$ws.Range("E8").Value = "  -3.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06608"
$ws.Range("E9").Value = "  -1.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.91"
$ws.Range("E10").Value = "  -4.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07632"
$ws.Range("E11").Value = "  -1.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.559"
$ws.Range("D13").Value = "1.735.71"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").Value = "1.949.87"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5766"
$ws.Range("E15").Value = "  -3.79%  "
$ws.Range("D16").Value = "0.0₅8162"
$ws.Range("E16").Value = "  -3.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.89"
$ws.Range("E17").Value = "  -2.28%  "
$ws.Range("D18").Value = "27.377.25"
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.42"
$ws.Range("E19").Value = "  -4.42%  "
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.662"
$ws.Range("E21").Value = "  -3.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.45"
$ws.Range("E22").Value = "  -4.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.980"
$ws.Range("E23").Value = "  -4.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.34"
$ws.Range("E25").Value = "  -3.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.732"
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("E27").Value = "  -2.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.282"
$ws.Range("E28").Value = "  -2.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.29"
$ws.Range("E29").Value = "  -5.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05402"
$ws.Range("E30").Value = "  -4.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.290"
$ws.Range("E31").Value = "  -1.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.489"
$ws.Range("E32").Value = "  -5.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.426"
$ws.Range("E33").Value = "  -2.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.642"
$ws.Range("E34").Value = "  -2.55%  "
$ws.Range("E35").Value = "  +0.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9487"
$ws.Range("E36").Value = "  -2.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.416"
$ws.Range("E37").Value = "  -1.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5839"
$ws.Range("E38").Value = "  -2.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01632"
$ws.Range("E39").Value = "  -2.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.860"
$ws.Range("E40").Value = "  -0.90%  "
$ws.Range("D41").Value = "1.044.94"
$ws.Range("E41").Value = "  -0.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.005"
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8418"
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.73"
$ws.Range("E44").Value = "  -1.37%  "
$ws.Range("D45").Value = "1.854.85"
$ws.Range("E45").Value = "  -1.41%  "
$ws.Range("E46").Value = "  +2.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.98"
$ws.Range("E47").Value = "  -2.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4522"
$ws.Range("E48").Value = "  +1.90%  "
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.073"
$ws.Range("E50").Value = "  -2.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05244"
$ws.Range("E51").Value = "  -1.57%  "
